{"js": "// Reorder the \"Requisitos\" bullet list in the LOM3092 document so that\n// the \"LOM3036 -  Propriedades Mec\u00e2nicas  (Requisito)\" line moves from\n// being the 3rd line to being the 1st line of that paragraph. The other\n// lines (LOM3011, LOM3013, LOM3046) keep their relative order.\n//\n// The four lines live inside ONE paragraph (style \"ListBullet\") as four\n// separate runs, each `<w:r><w:t>...</w:t><w:br/></w:r>`, separated by\n// manual line breaks (`<w:br/>`, which Office.js exposes as \"\\v\" in the\n// paragraph's .text). We must only reorder -- not merge -- those runs.\n\nconst MOVE_TEXT = \"LOM3036 -  Propriedades Mec\u00e2nicas  (Requisito)\";\nconst ANCHOR_TEXT = \"LOM3011 -  Ensaios Mec\u00e2nicos  (Requisito)\";\n\nconst body = context.document.body;\n\n// Locate the paragraph that holds the requirement bullet list (the one\n// that contains both the line we want to move and its new neighbour).\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paras.items.length; i++) {\n  const t = paras.items[i].text;\n  if (t.indexOf(MOVE_TEXT) !== -1 && t.indexOf(ANCHOR_TEXT) !== -1) {\n    target = paras.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the Requisitos paragraph containing '\" + MOVE_TEXT + \"'\");\n}\n\n// Split the paragraph's range on the manual line-break character so each\n// logical bullet line (text + trailing \"\\v\") is individually addressable.\nlet lineRanges = target.getTextRanges([\"\\v\"], false);\nlineRanges.load(\"items/text\");\nawait context.sync();\n\n// Find the exact line (run) whose text is the one we want to move.\nlet moveLineText = null;\nfor (let i = 0; i < lineRanges.items.length; i++) {\n  if (lineRanges.items[i].text.indexOf(MOVE_TEXT) === 0) {\n    moveLineText = lineRanges.items[i].text; // includes trailing \"\\v\"\n    break;\n  }\n}\nif (moveLineText === null) {\n  throw new Error(\"Could not locate the '\" + MOVE_TEXT + \"' line inside the paragraph\");\n}\n\n// Build a minimal single-run OOXML fragment carrying that same text/break\n// so the re-inserted copy keeps its own distinct <w:r> (matching the\n// original markup) instead of being merged into a neighbouring run.\nconst escaped = moveLineText\n  .replace(/&/g, \"&amp;\")\n  .replace(/</g, \"&lt;\")\n  .replace(/>/g, \"&gt;\")\n  .replace(\"\\v\", \"\"); // the break itself is re-added as <w:br/>, drop the literal char\n\nconst runOoxml =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body><w:p><w:r><w:t>' + escaped + '</w:t><w:br/></w:r></w:p></w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\n// Insert that run at the very start of the paragraph (before \"LOM3011 ...\").\nconst startRange = target.getRange(\"Start\");\nstartRange.insertOoxml(runOoxml, \"Before\");\nawait context.sync();\n\n// Re-split the (now mutated) paragraph and drop the original occurrence of\n// the moved line -- i.e. the *second* run whose text starts with MOVE_TEXT.\nconst paras2 = body.paragraphs;\nparas2.load(\"items/text\");\nawait context.sync();\n\nlet target2 = null;\nfor (let i = 0; i < paras2.items.length; i++) {\n  const t = paras2.items[i].text;\n  if (t.indexOf(MOVE_TEXT) !== -1 && t.indexOf(ANCHOR_TEXT) !== -1) {\n    target2 = paras2.items[i];\n    break;\n  }\n}\n\nlet lineRanges2 = target2.getTextRanges([\"\\v\"], false);\nlineRanges2.load(\"items/text\");\nawait context.sync();\n\nlet seen = 0;\nlet dup = null;\nfor (let i = 0; i < lineRanges2.items.length; i++) {\n  if (lineRanges2.items[i].text.indexOf(MOVE_TEXT) === 0) {\n    seen++;\n    if (seen === 2) {\n      dup = lineRanges2.items[i];\n      break;\n    }\n  }\n}\n\nif (!dup) {\n  throw new Error(\"Could not find the original (duplicate) '\" + MOVE_TEXT + \"' line to remove\");\n}\n\ndup.delete();\nawait context.sync();\n", "ps1": "# Reorder the \"Requisitos\" bullet list in the LOM3092 document so that\n# the \"LOM3036 -  Propriedades Mec\u00e2nicas  (Requisito)\" line moves from\n# being the 3rd line to being the 1st line of that paragraph. The other\n# lines (LOM3011, LOM3013, LOM3046) keep their relative order.\n#\n# The four lines live inside ONE paragraph (style \"ListBullet\") as four\n# separate runs, each <w:r><w:t>...</w:t><w:br/></w:r>, separated by\n# manual line breaks (w:br, which Word exposes as Chr(11)/vertical-tab in\n# Range.Text). We only reorder -- not merge -- those runs.\n\n$d = $word.ActiveDocument\n\n$moveText = \"LOM3036 -  Propriedades Mec\u00e2nicas  (Requisito)\"\n$anchorText = \"LOM3011 -  Ensaios Mec\u00e2nicos  (Requisito)\"\n$vt = [char]11\n\n# Locate the paragraph that holds the requirement bullet list (the one\n# containing both the line to move and its new neighbour).\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t.IndexOf($moveText) -ge 0 -and $t.IndexOf($anchorText) -ge 0) {\n        $target = $p\n        break\n    }\n}\nif ($target -eq $null) {\n    throw \"Could not find the Requisitos paragraph containing '$moveText'\"\n}\n\n$paraRange = $target.Range\n\n# Locate the \"LOM3036 ...\" run inside that paragraph with Find, then grow\n# the match by one character to also capture its trailing manual break.\n$searchRange = $paraRange.Duplicate\n$searchRange.Find.ClearFormatting()\n$searchRange.Find.Text = $moveText\n$found = $searchRange.Find.Execute()\nif (-not $found) {\n    throw \"Could not locate the '$moveText' line inside the paragraph\"\n}\n$lineRange = $d.Range($searchRange.Start, $searchRange.End + 1)\n$lineText = $lineRange.Text   # e.g. \"LOM3036 -  Propriedades Mec\u00e2nicas  (Requisito)\" + vt\n\n# Insert a copy of that text+break at the very start of the paragraph,\n# ahead of \"LOM3011 ...\". InsertBefore on a fresh Range keeps it as its\n# own run instead of merging into a neighbouring one.\n$paraStart = $d.Range($paraRange.Start, $paraRange.Start)\n$paraStart.InsertBefore($lineText)\n\n# Re-find the ORIGINAL \"LOM3036 ...\" occurrence (now pushed later in the\n# paragraph, after the copy we just inserted) and remove it so the line\n# exists only once, now at the front.\n$target2 = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t.IndexOf($moveText) -ge 0 -and $t.IndexOf($anchorText) -ge 0) {\n        $target2 = $p\n        break\n    }\n}\n$paraRange2 = $target2.Range\n\n$search2 = $paraRange2.Duplicate\n$search2.Find.ClearFormatting()\n$search2.Find.Text = $moveText\n# Skip past the newly-inserted copy so Find lands on the original.\n$search2.Start = $paraRange2.Start + $lineText.Length\n$search2.End = $paraRange2.End\n$found2 = $search2.Find.Execute()\nif (-not $found2) {\n    throw \"Could not find the original (duplicate) '$moveText' line to remove\"\n}\n\n$origLineRange = $d.Range($search2.Start, $search2.End + 1)\n$origLineRange.Delete()\n"}
